{"js": "// Update the p-values table in \"Fig 3\" (QUOTAS figures-with-above-advice update).\n// Table layout (row, col -- 0-indexed):\n//   row 0: header            \" \" | \"Cod\" | \"Hake\"\n//   row 1: \"GDP 2016\"        0.40 | 0.24\n//   row 2: \"OHI 2016\"        0.18 | 0.01\n//   row 3: \"OHI economic 2016\" 0.79 | 0.97\n//   row 4: \"Readiness\"       0.08 | 0.16\n//   row 5: \"Vulnerability\"   <0.01 | 0.01\n//\n// Target changes (per diff):\n//   (1,2) 0.24 -> 0.20\n//   (2,1) 0.18 -> 0.16\n//   (3,1) 0.79 -> 0.76\n//   (3,2) 0.97 -> 0.96\n//   (4,1) 0.08 -> 0.10\n//   (4,2) 0.16 -> 0.11\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Each edit is scoped to its own cell (row, col) so duplicate numeric\n// strings elsewhere in the table (e.g. \"0.16\" shows up both as a\n// pre-existing value and as a new value) can never cross-match.\nconst edits = [\n  { row: 1, col: 2, from: \"0.24\", to: \"0.20\" },\n  { row: 2, col: 1, from: \"0.18\", to: \"0.16\" },\n  { row: 3, col: 1, from: \"0.79\", to: \"0.76\" },\n  { row: 3, col: 2, from: \"0.97\", to: \"0.96\" },\n  { row: 4, col: 1, from: \"0.08\", to: \"0.10\" },\n  { row: 4, col: 2, from: \"0.16\", to: \"0.11\" },\n];\n\nfor (const edit of edits) {\n  const cell = table.getCell(edit.row, edit.col);\n  const results = cell.body.search(edit.from, {\n    matchCase: true,\n    matchWholeWord: true,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(\n      `Expected exactly 1 match for \"${edit.from}\" in cell (${edit.row},${edit.col}), found ${results.items.length}`\n    );\n  }\n\n  results.items[0].insertText(edit.to, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Update the p-values table in \"Fig 3\" (QUOTAS figures-with-above-advice update).\n# Table layout (Word COM 1-indexed: Table.Cell(row, col)):\n#   row 1: header               \" \" | \"Cod\" | \"Hake\"\n#   row 2: \"GDP 2016\"           0.40 | 0.24\n#   row 3: \"OHI 2016\"           0.18 | 0.01\n#   row 4: \"OHI economic 2016\"  0.79 | 0.97\n#   row 5: \"Readiness\"          0.08 | 0.16\n#   row 6: \"Vulnerability\"      <0.01 | 0.01\n#\n# Target changes (per diff):\n#   Cell(2,3) 0.24 -> 0.20\n#   Cell(3,2) 0.18 -> 0.16\n#   Cell(4,2) 0.79 -> 0.76\n#   Cell(4,3) 0.97 -> 0.96\n#   Cell(5,2) 0.08 -> 0.10\n#   Cell(5,3) 0.16 -> 0.11\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Each replacement is scoped to its own cell's Range so duplicate numeric\n# strings elsewhere in the table (e.g. \"0.16\" is both a pre-existing value\n# and a new value) can never cross-match.\n$edits = @(\n    @{ Row = 2; Col = 3; From = \"0.24\"; To = \"0.20\" },\n    @{ Row = 3; Col = 2; From = \"0.18\"; To = \"0.16\" },\n    @{ Row = 4; Col = 2; From = \"0.79\"; To = \"0.76\" },\n    @{ Row = 4; Col = 3; From = \"0.97\"; To = \"0.96\" },\n    @{ Row = 5; Col = 2; From = \"0.08\"; To = \"0.10\" },\n    @{ Row = 5; Col = 3; From = \"0.16\"; To = \"0.11\" }\n)\n\nforeach ($edit in $edits) {\n    $cell = $t.Cell($edit.Row, $edit.Col)\n    $rng = $cell.Range\n\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Text = $edit.From\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $edit.To\n    $find.Forward = $true\n    $find.Wrap = 0          # wdFindStop - stay within the cell range\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $true\n    $find.MatchWildcards = $false\n\n    $found = $find.Execute(\n        $edit.From,   # FindText\n        $false,       # MatchCase (kept false here; property above already set)\n        $true,        # MatchWholeWord\n        $false,       # MatchWildcards\n        $false,       # MatchSoundsLike\n        $false,       # MatchAllWordForms\n        $true,        # Forward\n        0,            # Wrap = wdFindStop\n        $false,       # Format\n        $edit.To,     # ReplaceWith\n        2             # Replace = wdReplaceAll\n    )\n\n    if (-not $found) {\n        throw \"Could not find '$($edit.From)' in cell ($($edit.Row),$($edit.Col))\"\n    }\n}\n"}
